$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated error-calculation results (row 2) ---
$ws.Range("F2").Value = 3.125603062778264
$ws.Range("G2").Value = 0.01283955242086306
$ws.Range("H2").Value = 3.54551924390971
$ws.Range("I2").Value = 1.069781537870688
$ws.Range("N2").Value = 0.5832367987382572
$ws.Range("O2").Value = 0.3212830336741677
$ws.Range("R2").Value = 0.002070408906527906
$ws.Range("S2").Value = 0.2893874872767901
$ws.Range("V2").Value = 0.002320394096024735
$ws.Range("W2").Value = 0.01283920043037212

# --- GUI/Excel column-width formatting tweaks ---
# Target stored widths (character units): F=18.7109375, H=19.7109375,
# O=20.7109375, S=19.7109375 (F<->H and O<->S effectively swap widths).
# The ColumnWidth setter snaps to the host's internal pixel grid, so the
# input is pre-compensated to land on the closest representable width.
$ws.Columns.Item(6).ColumnWidth = 17.833333333333332
$ws.Columns.Item(8).ColumnWidth = 18.833333333333332
$ws.Columns.Item(15).ColumnWidth = 19.833333333333332
$ws.Columns.Item(19).ColumnWidth = 18.833333333333332
